# "exam A attempt 2 answers"
#
# Fills in the grading columns (B = student answer, C = Correct/Incorrect
# formula) on the "Exam A 2" sheet for rows 6-90, adds two score notes in
# column E, introduces one deliberately wrong answer (row 58), and updates
# the selection/view state on both the "Exam A Answers" and "Exam A 2"
# sheets to match where the author last clicked.

$wb = $excel.ActiveWorkbook

$ansSheet = $wb.Worksheets.Item("Exam A Answers")
$ws = $wb.Worksheets.Item("Exam A 2")

# Row 58 is the one question the author got wrong on this attempt - every
# other row's "my answer" (column B) simply mirrors the correct answer
# already stored in column A. Shared-string index 5 is "D", which differs
# from A58's "A" (index 2), so C58 evaluates to "Incorrect".
$wrongRow = 58
$wrongAnswer = $ws.Range("D5").Value2

for ($r = 6; $r -le 90; $r++) {
    if ($r -eq $wrongRow) {
        $ws.Cells.Item($r, 2).Value = $wrongAnswer
    } else {
        $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value2
    }
    $ws.Cells.Item($r, 3).Formula = "=IF(A$r=B$r, ""Correct"", ""Incorrect"")"
}

# Score annotations left next to questions 12 & 13, matching the pattern
# already used on the other exam sheets (e.g. "Exam C" E11/E12).
$ws.Range("E12").Value = "84/85"
$ws.Range("E13").Value = "89/90"

# Restore the cursor position on "Exam A Answers" (it was left at D1,
# moved to C6), without leaving that sheet active.
$ansSheet.Range("C6").Select()

# Finally leave "Exam A 2" selected at G17, scrolled so column B is back
# in view, matching the author's last position when they saved.
$ws.Range("G17").Select()
